# The workbook contains a weekly price log for "Pepino dulce" at the
# Vega Monumental Concepción market. A new week's record was inserted as
# row 4 (pushing the previous rows 4-74 down to rows 5-75), and the new
# row was populated with that week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4, shifting existing rows 4-74 down to 5-75.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new week's data.
$ws.Cells.Item(4, 1).Value2 = 11
$ws.Cells.Item(4, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(4, 3).Value2 = "Bíobío"
$ws.Cells.Item(4, 4).Value2 = 45083
$ws.Cells.Item(4, 5).Value2 = 8
$ws.Cells.Item(4, 6).Value2 = 100112043
$ws.Cells.Item(4, 7).Value2 = "Pepino dulce"
$ws.Cells.Item(4, 8).Value2 = "Cultivar IV Región"
$ws.Cells.Item(4, 9).Value2 = "Primera"
$ws.Cells.Item(4, 10).Value2 = 140
$ws.Cells.Item(4, 11).Value2 = 12000
$ws.Cells.Item(4, 12).Value2 = 14000
$ws.Cells.Item(4, 13).Value2 = 13143
$ws.Cells.Item(4, 14).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(4, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(4, 16).Value2 = 730
$ws.Cells.Item(4, 17).Value2 = 18
$ws.Cells.Item(4, 18).Value2 = "Hortaliza"
